$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serping1"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.592460999999999
$ws.Range("H2").Value = 13.777383
$ws.Range("I2").Value = 0.003302946473568516
$ws.Range("J2").Value = 0.003302946473568516
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.866432
$ws.Range("N2").Value = 8.599295999999999
$ws.Range("O2").Value = 0.9456981836489474
$ws.Range("P2").Value = 0.9456981836489475
$ws.Range("Q2").Value = 13.163977169152
$ws.Range("R2").Value = 118.475794522368
$ws.Range("S2").Value = 0.003123590480743441
$ws.Range("T2").Value = 0.003123590480743442

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serping1"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.592460999999999
$ws.Range("H3").Value = 13.777383
$ws.Range("I3").Value = 0.003302946473568516
$ws.Range("J3").Value = 0.003302946473568516
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.16459
$ws.Range("N3").Value = 0.49377
$ws.Range("O3").Value = 0.05430181635105255
$ws.Range("P3").Value = 0.05430181635105256
$ws.Range("Q3").Value = 0.7558731559899998
$ws.Range("R3").Value = 6.802858403909999
$ws.Range("S3").Value = 0.0001793559928250742
$ws.Range("T3").Value = 0.0001793559928250742

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Serping1"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1205.102620666667
$ws.Range("H4").Value = 3615.307862
$ws.Range("I4").Value = 0.8667225374846176
$ws.Range("J4").Value = 0.8667225374846176
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.866432
$ws.Range("N4").Value = 8.599295999999999
$ws.Range("O4").Value = 0.9456981836489474
$ws.Range("P4").Value = 0.9456981836489475
$ws.Range("Q4").Value = 3454.344715162795
$ws.Range("R4").Value = 31089.10243646515
$ws.Range("S4").Value = 0.8196579294268096
$ws.Range("T4").Value = 0.8196579294268097

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serping1"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1205.102620666667
$ws.Range("H5").Value = 3615.307862
$ws.Range("I5").Value = 0.8667225374846176
$ws.Range("J5").Value = 0.8667225374846176
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.16459
$ws.Range("N5").Value = 0.49377
$ws.Range("O5").Value = 0.05430181635105255
$ws.Range("P5").Value = 0.05430181635105256
$ws.Range("Q5").Value = 198.3478403355267
$ws.Range("R5").Value = 1785.13056301974
$ws.Range("S5").Value = 0.04706460805780797
$ws.Range("T5").Value = 0.04706460805780797

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Serping1"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 180.7183073333333
$ws.Range("H6").Value = 542.154922
$ws.Range("I6").Value = 0.1299745160418139
$ws.Range("J6").Value = 0.1299745160418139
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.866432
$ws.Range("N6").Value = 8.599295999999999
$ws.Range("O6").Value = 0.9456981836489474
$ws.Range("P6").Value = 0.9456981836489475
$ws.Range("Q6").Value = 518.0167391261012
$ws.Range("R6").Value = 4662.150652134912
$ws.Range("S6").Value = 0.1229166637413944
$ws.Range("T6").Value = 0.1229166637413944

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Serping1"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 180.7183073333333
$ws.Range("H7").Value = 542.154922
$ws.Range("I7").Value = 0.1299745160418139
$ws.Range("J7").Value = 0.1299745160418139
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.16459
$ws.Range("N7").Value = 0.49377
$ws.Range("O7").Value = 0.05430181635105255
$ws.Range("P7").Value = 0.05430181635105256
$ws.Range("Q7").Value = 29.74442620399333
$ws.Range("R7").Value = 267.69983583594
$ws.Range("S7").Value = 0.007057852300419512
$ws.Range("T7").Value = 0.007057852300419515

# Remove old rows 8-10 (data no longer present)
$ws.Range("A8:T10").Delete()
